# Auto-generated edit script: updates Leve profit-tracking values across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets to reflect refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1067.9286
$ws.Range("I18").Value = 906.5
$ws.Range("J18").Value = 1283.1666
$ws.Range("K18").Value = 906.5
$ws.Range("L18").Value = 1283.1666
$ws.Range("M18").Value = -622.5
$ws.Range("N18").Value = -1851.1666
$ws.Range("H43").Value = 1324.8572
$ws.Range("I43").Value = 1315
$ws.Range("K43").Value = 1315
$ws.Range("M43").Value = -1246
$ws.Range("H70").Value = 1949.2
$ws.Range("I70").Value = 2124
$ws.Range("K70").Value = 6372
$ws.Range("M70").Value = -6102
$ws.Range("H73").Value = 1949.2
$ws.Range("I73").Value = 2124
$ws.Range("K73").Value = 6372
$ws.Range("M73").Value = -5436
$ws.Range("H80").Value = 1548.75
$ws.Range("I80").Value = 347.5
$ws.Range("K80").Value = 1042.5
$ws.Range("M80").Value = -44.5
$ws.Range("H83").Value = 1548.75
$ws.Range("I83").Value = 347.5
$ws.Range("K83").Value = 3127.5
$ws.Range("M83").Value = 1864.5
$ws.Range("H88").Value = 1909.3334
$ws.Range("J88").Value = 1909.3334
$ws.Range("L88").Value = 1909.3334
$ws.Range("N88").Value = -2721.3334
$ws.Range("H91").Value = 1909.3334
$ws.Range("J91").Value = 1909.3334
$ws.Range("L91").Value = 1909.3334
$ws.Range("N91").Value = -4717.3334
$ws.Range("H94").Value = 1446.5
$ws.Range("I94").Value = 1446.5
$ws.Range("K94").Value = 1446.5
$ws.Range("M94").Value = -995.5
$ws.Range("H103").Value = 500
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 1500
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -2672
$ws.Range("H125").Value = 4000
$ws.Range("I125").Value = 4000
$ws.Range("K125").Value = 36000
$ws.Range("M125").Value = -33540
$ws.Range("H141").Value = 2999
$ws.Range("I141").Value = 2999
$ws.Range("K141").Value = 8997
$ws.Range("M141").Value = -3817

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2226.2942
$ws.Range("I2").Value = 1684.2
$ws.Range("J2").Value = 3000.7144
$ws.Range("K2").Value = 1684.2
$ws.Range("L2").Value = 3000.7144
$ws.Range("M2").Value = -1571.2
$ws.Range("N2").Value = -3226.7144
$ws.Range("H5").Value = 129.14285
$ws.Range("I5").Value = 17.333334
$ws.Range("K5").Value = 17.333334
$ws.Range("M5").Value = 94.66666599999999
$ws.Range("H97").Value = 800.3333
$ws.Range("I97").Value = 751.4286
$ws.Range("K97").Value = 751.4286
$ws.Range("M97").Value = -255.4286
$ws.Range("H116").Value = 2226.2942
$ws.Range("I116").Value = 1684.2
$ws.Range("J116").Value = 3000.7144
$ws.Range("K116").Value = 1684.2
$ws.Range("L116").Value = 3000.7144
$ws.Range("M116").Value = 609.8
$ws.Range("N116").Value = -7588.7144
$ws.Range("H122").Value = 1776.6
$ws.Range("I122").Value = 1095.875
$ws.Range("K122").Value = 3287.625
$ws.Range("M122").Value = -837.625
$ws.Range("H132").Value = 850
$ws.Range("I132").Value = 850
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2550
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -20
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2226.2942
$ws.Range("I3").Value = 1684.2
$ws.Range("J3").Value = 3000.7144
$ws.Range("K3").Value = 1684.2
$ws.Range("L3").Value = 3000.7144
$ws.Range("M3").Value = -1570.2
$ws.Range("N3").Value = -3228.7144
$ws.Range("H4").Value = 129.14285
$ws.Range("I4").Value = 17.333334
$ws.Range("K4").Value = 17.333334
$ws.Range("M4").Value = 97.66666599999999
$ws.Range("H86").Value = 1182.3684
$ws.Range("I86").Value = 974.9167
$ws.Range("K86").Value = 974.9167
$ws.Range("M86").Value = 148.0833
$ws.Range("H89").Value = 1182.3684
$ws.Range("I89").Value = 974.9167
$ws.Range("K89").Value = 4874.5835
$ws.Range("M89").Value = 741.4165000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1767.8182
$ws.Range("I7").Value = 492.57144
$ws.Range("K7").Value = 492.57144
$ws.Range("M7").Value = -379.57144
$ws.Range("H22").Value = 889.4286
$ws.Range("I22").Value = 826
$ws.Range("K22").Value = 826
$ws.Range("M22").Value = -476
$ws.Range("H31").Value = 2008.8334
$ws.Range("I31").Value = 1910.6
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1910.6
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1615.6
$ws.Range("N31").Value = -3090
$ws.Range("H34").Value = 2008.8334
$ws.Range("I34").Value = 1910.6
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1910.6
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1708.6
$ws.Range("N34").Value = -2904
$ws.Range("H122").Value = 2386.8
$ws.Range("I122").Value = 2386.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7160.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4710.400000000001
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 220
$ws.Range("I33").Value = 220
$ws.Range("K33").Value = 1320
$ws.Range("M33").Value = -1037
$ws.Range("H120").Value = 3333
$ws.Range("I120").Value = 3333
$ws.Range("K120").Value = 9999
$ws.Range("M120").Value = -5161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 94.888885
$ws.Range("I2").Value = 80.25
$ws.Range("J2").Value = 106.6
$ws.Range("K2").Value = 80.25
$ws.Range("L2").Value = 106.6
$ws.Range("M2").Value = 32.75
$ws.Range("N2").Value = -332.6
$ws.Range("H11").Value = 2762500.8
$ws.Range("J11").Value = 514500
$ws.Range("L11").Value = 514500
$ws.Range("N11").Value = -514778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7757.522
$ws.Range("I122").Value = 8135
$ws.Range("K122").Value = 24405
$ws.Range("M122").Value = -21955
$ws.Range("H132").Value = 3235.353
$ws.Range("I132").Value = 3062.077
$ws.Range("K132").Value = 9186.231
$ws.Range("M132").Value = -6656.231
